$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# The page with "Requisitos" ends with a page-break paragraph, an empty
# paragraph, a "Ver no Jupiter..." paragraph and a copyright paragraph.
# The footer text (empty paragraph + both text paragraphs) is removed,
# leaving only the page-break paragraph in place.
$jupiterPara = Find-ParagraphByText $d "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightPara = Find-ParagraphByText $d "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
$emptyPara = $d.Paragraphs.Item($jupiterPara.Index - 1)

# Delete from the last paragraph upward so earlier object references /
# indices remain valid as each delete happens.
$copyrightPara.Range.Delete()
$jupiterPara.Range.Delete()
$emptyPara.Range.Delete()

Write-Output ("Paragraphs remaining: " + $d.Paragraphs.Count)
